$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# Overview row 4: 2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md
$wsOverview.Range("A4").Value = "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md", "", "", "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md")
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-41-13 06:41:02"

# zh-cn row 4: 2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md
$wsZh.Range("A4").Value = "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md", "", "", "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md")
$wsZh.Range("B4").Value = ".md"
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md", "", "", ".md")
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.19da5d73730fcf81df93bbc4a1e8aa24cc6aedf1.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.19da5d73730fcf81df93bbc4a1e8aa24cc6aedf1.zh-cn.xlf", "", "", "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.19da5d73730fcf81df93bbc4a1e8aa24cc6aedf1.zh-cn.xlf")
$wsZh.Range("E4").Value = "2016-03-13 06:40:58"
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"

# de-de row 4: 2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md
$wsDe.Range("A4").Value = "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md", "", "", "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md")
$wsDe.Range("B4").Value = ".md"
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.md", "", "", ".md")
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.19da5d73730fcf81df93bbc4a1e8aa24cc6aedf1.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.19da5d73730fcf81df93bbc4a1e8aa24cc6aedf1.de-de.xlf", "", "", "2cca0b2b-b8fd-4a61-820d-e6b34c7054aa.19da5d73730fcf81df93bbc4a1e8aa24cc6aedf1.de-de.xlf")
$wsDe.Range("E4").Value = "2016-03-13 06:41:02"
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"

# Overview row 5: 7e37c851-2e47-493a-aa87-a829ea808040.md
$wsOverview.Range("A5").Value = "7e37c851-2e47-493a-aa87-a829ea808040.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/7e37c851-2e47-493a-aa87-a829ea808040.md", "", "", "7e37c851-2e47-493a-aa87-a829ea808040.md")
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-41-13 06:41:02"

# zh-cn row 5: 7e37c851-2e47-493a-aa87-a829ea808040.md
$wsZh.Range("A5").Value = "7e37c851-2e47-493a-aa87-a829ea808040.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/7e37c851-2e47-493a-aa87-a829ea808040.md", "", "", "7e37c851-2e47-493a-aa87-a829ea808040.md")
$wsZh.Range("B5").Value = ".md"
$wsZh.Hyperlinks.Add($wsZh.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/7e37c851-2e47-493a-aa87-a829ea808040.md", "", "", ".md")
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "7e37c851-2e47-493a-aa87-a829ea808040.21a0ae5e58bbf979c1679ec4f52ad4ceadb9d84b.zh-cn.xlf"
$wsZh.Hyperlinks.Add($wsZh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7e37c851-2e47-493a-aa87-a829ea808040.21a0ae5e58bbf979c1679ec4f52ad4ceadb9d84b.zh-cn.xlf", "", "", "7e37c851-2e47-493a-aa87-a829ea808040.21a0ae5e58bbf979c1679ec4f52ad4ceadb9d84b.zh-cn.xlf")
$wsZh.Range("E5").Value = "2016-03-13 06:40:58"
$wsZh.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H5").Value = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value = "Include"

# de-de row 5: 7e37c851-2e47-493a-aa87-a829ea808040.md
$wsDe.Range("A5").Value = "7e37c851-2e47-493a-aa87-a829ea808040.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/7e37c851-2e47-493a-aa87-a829ea808040.md", "", "", "7e37c851-2e47-493a-aa87-a829ea808040.md")
$wsDe.Range("B5").Value = ".md"
$wsDe.Hyperlinks.Add($wsDe.Range("B5"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/7e37c851-2e47-493a-aa87-a829ea808040.md", "", "", ".md")
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "7e37c851-2e47-493a-aa87-a829ea808040.21a0ae5e58bbf979c1679ec4f52ad4ceadb9d84b.de-de.xlf"
$wsDe.Hyperlinks.Add($wsDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7e37c851-2e47-493a-aa87-a829ea808040.21a0ae5e58bbf979c1679ec4f52ad4ceadb9d84b.de-de.xlf", "", "", "7e37c851-2e47-493a-aa87-a829ea808040.21a0ae5e58bbf979c1679ec4f52ad4ceadb9d84b.de-de.xlf")
$wsDe.Range("E5").Value = "2016-03-13 06:41:02"
$wsDe.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value = "Include"

